# Update "想去人数" (F column) counts on gh-pages data refresh
# Sheet "展览" (Exhibition)
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 127
$ws1.Range("F3").Value = 190
$ws1.Range("F4").Value = 440
$ws1.Range("F7").Value = 1228
$ws1.Range("F8").Value = 412
$ws1.Range("F9").Value = 205
$ws1.Range("F12").Value = 384
$ws1.Range("F13").Value = 428
$ws1.Range("F15").Value = 188
$ws1.Range("F16").Value = 735
$ws1.Range("F20").Value = 479
$ws1.Range("F21").Value = 278
$ws1.Range("F22").Value = 93
$ws1.Range("F26").Value = 484
$ws1.Range("F27").Value = 9

# Sheet "演出" (Performance)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F5").Value = 44
$ws2.Range("F6").Value = 45
$ws2.Range("F10").Value = 632
$ws2.Range("F11").Value = 152
$ws2.Range("F12").Value = 57

# Sheet "全部类型" (All Types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 127
$ws4.Range("F5").Value = 190
$ws4.Range("F6").Value = 440
$ws4.Range("F9").Value = 1228
$ws4.Range("F10").Value = 412
$ws4.Range("F11").Value = 205
$ws4.Range("F16").Value = 44
$ws4.Range("F17").Value = 384
$ws4.Range("F18").Value = 45
$ws4.Range("F20").Value = 428
$ws4.Range("F22").Value = 188
$ws4.Range("F23").Value = 735
$ws4.Range("F27").Value = 479
$ws4.Range("F30").Value = 278
$ws4.Range("F31").Value = 93
$ws4.Range("F33").Value = 632
$ws4.Range("F34").Value = 152
$ws4.Range("F37").Value = 57
$ws4.Range("F39").Value = 484
$ws4.Range("F42").Value = 9
